# Applies the "Updated cryptos list" data refresh (GitHub Actions run).
# For every changed cell we:
#   1. force the cell's number format to Text ("@") so a numeric-looking
#      string (e.g. "1.00", "72.30", "0.0697") is stored verbatim instead
#      of being parsed/collapsed into a Number by Excel's input parser,
#   2. assign the new literal value,
#   3. reapply the "Normal" cell style so we don't leave a stray
#      Text-format style behind on cells that started out unstyled
#      (keeps styles identical to the original file).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue([string]$cellRef, [string]$value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" "66.611.44"
Set-TextValue "E2" "  -0.60%  "
Set-TextValue "D3" "3.450.14"
Set-TextValue "E3" "  -0.25%  "
Set-TextValue "E4" "  +0.05%  "
Set-TextValue "D5" "584.78"
Set-TextValue "E5" "  -1.45%  "
Set-TextValue "D6" "176.54"
Set-TextValue "E6" "  -1.91%  "
Set-TextValue "E7" "  +0.01%  "
Set-TextValue "D8" "0.602"
Set-TextValue "E8" "  -1.02%  "
Set-TextValue "D9" "3.449.45"
Set-TextValue "E9" "  -0.18%  "
Set-TextValue "E10" "  -6.03%  "
Set-TextValue "D11" "6.88"
Set-TextValue "E11" "  -1.10%  "
Set-TextValue "E12" "  -3.15%  "
Set-TextValue "D13" "4.046.01"
Set-TextValue "E13" "  -0.23%  "
Set-TextValue "E14" "  -4.63%  "
Set-TextValue "E15" "  -0.29%  "
Set-TextValue "D16" "66.579.85"
Set-TextValue "E16" "  -0.60%  "
Set-TextValue "E17" "  -2.47%  "
Set-TextValue "D18" "3.447.21"
Set-TextValue "E18" "  -0.22%  "
Set-TextValue "D19" "5.97"
Set-TextValue "E19" "  -4.10%  "
Set-TextValue "E20" "  -2.43%  "
Set-TextValue "D21" "378.81"
Set-TextValue "E21" "  -2.77%  "
Set-TextValue "D22" "7.83"
Set-TextValue "E22" "  -1.11%  "
Set-TextValue "D23" "1.00"
Set-TextValue "E23" "  +0.14%  "
Set-TextValue "B24" "Litecoin"
Set-TextValue "C24" "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue "D24" "72.30"
Set-TextValue "E24" "  +0.38%  "
Set-TextValue "B25" "LEO"
Set-TextValue "C25" "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue "D25" "5.75"
Set-TextValue "E25" "  -0.10%  "
Set-TextValue "B26" "Polygon"
Set-TextValue "C26" "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextValue "D26" "0.537"
Set-TextValue "E26" "  +0.17%  "
Set-TextValue "E27" "  -1.26%  "
Set-TextValue "D28" "9.80"
Set-TextValue "E28" "  -5.15%  "
Set-TextValue "E29" "  +0.82%  "
Set-TextValue "E30" "  +0.31%  "
Set-TextValue "D31" "24.32"
Set-TextValue "E31" "  +3.99%  "
Set-TextValue "D32" "5.86"
Set-TextValue "E32" "  -4.53%  "
Set-TextValue "D33" "2.00"
Set-TextValue "E33" "  -2.76%  "
Set-TextValue "E34" "  -5.60%  "
Set-TextValue "E35" "  +0.00%  "
Set-TextValue "E36" "  -1.68%  "
Set-TextValue "E37" "  -0.21%  "
Set-TextValue "D38" "161.17"
Set-TextValue "E38" "  -1.26%  "
Set-TextValue "D39" "29.52"
Set-TextValue "E39" "  +13.08%  "
Set-TextValue "D40" "0.892"
Set-TextValue "E40" "  +2.07%  "
Set-TextValue "E42" "  -6.21%  "
Set-TextValue "D43" "4.52"
Set-TextValue "E43" "  -2.70%  "
Set-TextValue "D44" "2.737.03"
Set-TextValue "E44" "  -0.18%  "
Set-TextValue "D45" "6.44"
Set-TextValue "E45" "  -5.02%  "
Set-TextValue "D46" "0.0697"
Set-TextValue "E46" "  -3.16%  "
Set-TextValue "D47" "40.67"
Set-TextValue "E47" "  -1.50%  "
Set-TextValue "D48" "24.55"
Set-TextValue "E48" "  -6.27%  "
Set-TextValue "D49" "0.0294"
Set-TextValue "E49" "  -1.51%  "
Set-TextValue "D50" "309.01"
Set-TextValue "E50" "  -5.65%  "
Set-TextValue "D51" "0.830"
Set-TextValue "E51" "  -1.11%  "
